$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows before row 292 (shifts existing 292-298 down to 294-300)
$ws.Rows.Item(292).EntireRow.Insert()
$ws.Rows.Item(292).EntireRow.Insert()

# Fill row 292
$ws.Cells.Item(292,1).Value2 = 1
$ws.Cells.Item(292,2).Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(292,3).Value2 = "Arica y Parinacota"
$ws.Cells.Item(292,4).Value2 = 44595
$ws.Cells.Item(292,5).Value2 = 15
$ws.Cells.Item(292,6).Value2 = 100112023
$ws.Cells.Item(292,7).Value2 = "Brócoli"
$ws.Cells.Item(292,8).Value2 = "Sin especificar"
$ws.Cells.Item(292,9).Value2 = "Segunda"
$ws.Cells.Item(292,10).Value2 = 1200
$ws.Cells.Item(292,11).Value2 = 700
$ws.Cells.Item(292,12).Value2 = 800
$ws.Cells.Item(292,13).Value2 = 750
$ws.Cells.Item(292,14).Value2 = "$/unidad"
$ws.Cells.Item(292,15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(292,16).Value2 = 750
$ws.Cells.Item(292,17).Value2 = 1
$ws.Cells.Item(292,18).Value2 = "Hortaliza"

# Fill row 293
$ws.Cells.Item(293,1).Value2 = 1
$ws.Cells.Item(293,2).Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(293,3).Value2 = "Arica y Parinacota"
$ws.Cells.Item(293,4).Value2 = 44595
$ws.Cells.Item(293,5).Value2 = 15
$ws.Cells.Item(293,6).Value2 = 100112023
$ws.Cells.Item(293,7).Value2 = "Brócoli"
$ws.Cells.Item(293,8).Value2 = "Sin especificar"
$ws.Cells.Item(293,9).Value2 = "Tercera"
$ws.Cells.Item(293,10).Value2 = 800
$ws.Cells.Item(293,11).Value2 = 450
$ws.Cells.Item(293,12).Value2 = 500
$ws.Cells.Item(293,13).Value2 = 475
$ws.Cells.Item(293,14).Value2 = "$/unidad"
$ws.Cells.Item(293,15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(293,16).Value2 = 475
$ws.Cells.Item(293,17).Value2 = 1
$ws.Cells.Item(293,18).Value2 = "Hortaliza"

Write-Host "done"
Write-Host "rows now:" $ws.UsedRange.Rows.Count
